$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix tiny floating point rounding on existing row 15 timestamp
$ws.Range("A15").Value = 45863.9169284375

# Append new row 16 with the latest sensor reading
$ws.Range("A16").Value = 45863.95854514534
$ws.Range("B16").Value = 2025
$ws.Range("C16").Value = 30
$ws.Range("D16").Value = 13.22
$ws.Range("E16").Value = 89.66
$ws.Range("F16").Value = 0
$ws.Range("G16").Value = 1.91
$ws.Range("H16").Value = "NNW"
$ws.Range("I16").Value = 0
$ws.Range("J16").Value = "23:00:18"

# Row 16's date cell should use the same date/time number format as the rest of column A
$ws.Range("A16").NumberFormat = "YYYY-MM-DD HH:MM:SS"
